$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.891.43'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '3.141.76'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.83'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.14'
$ws.Range('E6').Value = '  -1.68%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '3.133.26'
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.87'
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.20'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').Value = '3.661.72'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.34'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = '63.742.41'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '3.137.92'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '468.34'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.36'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.52'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '81.58'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.96'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('E26').Value = '  +6.81%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.70'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.42'
$ws.Range('E29').Value = '  +8.57%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.23'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.81'
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').Value = '0.0₃0847'
$ws.Range('E35').Value = '  -4.25%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +1.10%  '
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.21'
$ws.Range('E39').Value = '  -5.84%  '
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.32'
$ws.Range('E41').Value = '  +7.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '453.20'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('E43').Value = '  +5.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0371'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').Value = '2.908.70'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '39.55'
$ws.Range('E46').Value = '  +9.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.107'
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.71'
$ws.Range('E48').Value = '  +5.88%  '
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('E51').Value = '  -0.86%  '

Write-Output "Applied 83 cell updates"
